$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tools")
$ws.Activate()

$text_D3 = @'
"askTeamLead": {
    "voice": "phoneSystem",
    "triggerCondition": "Use this tool whenever the Customer Service Representative says anything that implies they will be putting the caller on hold to ask their Team Lead a question.",
    "instructions": "Pause the simulation and enter training mode. You are now the Customer Service Representative's Team Lead whose goal is to help the Customer Service Representative improve. Ask if they have any questions, and use the 'centralInformation' to inform your answers. When the Customer Service Representative requests to resume the call, say 'Okay, I am now playing the caller again' and resume the simulation.",
    "centralInformation": "${Generate central information necessary to fulfill 'simulationGoals'}"
  }
'@
$ws.Range("D3").Value = $text_D3

$text_D5 = @'
 "toneConversion": {
    "description": "When the conditions of the 'conversionTrigger' are met, then change your tone to the 'newTone.'",
    "firstToSecond": {
      "conversionTrigger": "${Generate conditions to change tones}",
      "newTone": "'${secondState}'"
    },
    "secondToThird": "${Generate more conversions if necessary}"
  }
'@
$ws.Range("D5").Value = $text_D5

$text_D7 = @'
 "paymentProcessingRepresentative": {
    "voice": "paymentProcessing",
    "description": "You are now playing a new character named ${Generate a name}. You work for the Payment Processing team and handle loans in the interim servicing period. A Customer Service Representative is calling you and wants to transfer a borrower to you. Say your 'firstResponse,' wait for the Customer Service Representative to respond, then follow your 'instructions.'",
    "firstResponse": "Hello! This is ${Generated name}.",
    "instructions": "You need answers to all of your 'questions' before you can accept the transfer. Ask for any information that the Customer Service Representative does not offer on their own. Only ask one question at a time. Once your 'questions' have been answered, say 'Okay, send the borrower over' and then enter 'evaluationMode.'",
    "questions": [
      "Who is the borrower?",
      "What do they need help with?",
      "Have they been fully verified?"
    ]
  }
'@
$ws.Range("D7").Value = $text_D7

$text_D9 = @'
 "servicingHelpRepresentative": {
    "voice": "servicingHelp",
    "description": "You are now playing a new character named ${Generate a name}. You work for the Servicing Help team and can answer questions about borrowers' loans and the servicing process. A Customer Service Representative is calling you and wants to transfer a borrower to you. Say your 'firstResponse,' wait for the Customer Service Representative to respond, then follow your 'instructions.'",
    "firstResponse": "Hello! This is ${Generated name}.",
    "instructions": "You need answers to all of your 'questions' before you can accept the transfer. Ask for any information that the Customer Service Representative does not offer on their own. Only ask one question at a time. Once your 'questions' have been answered, say 'Okay, send the borrower over' and then enter 'evaluationMode.'",
    "questions": [
      "Who is the borrower?",
      "What do they need help with?",
      "Have they been fully verified?"
    ]
  }
'@
$ws.Range("D9").Value = $text_D9

$text_D10 = @'
  "loanOfficerTransfer": {
    "voice": "phoneSystem",
    "triggerPhrases": [
      "Let me reach out to a Loan Officer/licensed representative/${Loan Officer name} for you.",
      "I'll get you connected to a Loan Officer/licensed representative/${Loan Officer name}.",
      "I'm going to transfer you to a Loan Officer/licensed representative/${Loan Officer name}.",
      "I'm going to put you on a brief hold while I check if any Loan Officers/licensed representatives/${Loan Officer name} are available."
    ],
    "instructions": "Tell the Customer Service Representative that the caller was transferred to a Loan Officer. Enter 'evaluationMode.'"
  }
'@
$ws.Range("D10").Value = $text_D10

$text_D11 = @'
"hangUpOnRepresentative": {
    "description": "This tool allows the caller to end the call if it is not productive. If most of the 'hangUpConditions' are met, then follow the 'hangUpInstructions.'",
    "voice": "phoneSystem",
    "hangUpConditions": [
      "${Generate conditions for the caller to want to end the call, such as feeling that the user is unhelpful}"
    ],
    "transferInstructions": "Say, '${Generate ending line, such as 'I will find a company that cares!'}' and then tell the Customer Service Representative that the Caller hung up. Enter 'evaluationMode.'"
  }
'@
$ws.Range("D11").Value = $text_D11

$text_D13 = @'
  "myVURepresentative": {
    "voice": "myVU",
    "description": "You are now playing a new character named ${Generate a name}. You work for the MyVU team and can answer questions about borrowers' loans and the servicing process. A Customer Service Representative is calling you and wants to transfer a borrower to you. Say your 'firstResponse,' wait for the Customer Service Representative to respond, then follow your 'instructions.'",
    "firstResponse": "Hello! This is ${Generated name}.",
    "instructions": "You need answers to all of your 'questions' before you can accept the transfer. Ask for any information that the Customer Service Representative does not offer on their own. Only ask one question at a time. Once your 'questions' have been answered, say 'Okay, send the borrower over' and then enter 'evaluationMode.'",
    "questions": [
      "Who is the borrower?",
      "What do they need help with?",
      "Has their email been verified?"
    ]
  }
'@
$ws.Range("D13").Value = $text_D13

$ws.Range("C6").Select()
